$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "26.044.31"
$ws.Range("E2").NumberFormat = "@"
$ws.Range("E2").Value = "  -0.16%  "

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "1.637.04"
$ws.Range("E3").NumberFormat = "@"
$ws.Range("E3").Value = "  -1.76%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "213.74"
$ws.Range("E5").NumberFormat = "@"
$ws.Range("E5").Value = "  +1.99%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "0.5243"
$ws.Range("E6").NumberFormat = "@"
$ws.Range("E6").Value = "  -0.31%  "

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "1.001"
$ws.Range("E7").NumberFormat = "@"
$ws.Range("E7").Value = "  -0.13%  "

$ws.Range("E8").NumberFormat = "@"
$ws.Range("E8").Value = "  -1.07%  "

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.06300"
$ws.Range("E9").NumberFormat = "@"
$ws.Range("E9").Value = "  +0.21%  "

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "20.66"
$ws.Range("E10").NumberFormat = "@"
$ws.Range("E10").Value = "  -2.39%  "

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.07656"
$ws.Range("E11").NumberFormat = "@"
$ws.Range("E11").Value = "  +1.56%  "

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "1.634.15"
$ws.Range("E12").NumberFormat = "@"
$ws.Range("E12").Value = "  -2.05%  "

$ws.Range("E13").NumberFormat = "@"
$ws.Range("E13").Value = "  -0.48%  "

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "1.860.35"
$ws.Range("E14").NumberFormat = "@"
$ws.Range("E14").Value = "  -1.81%  "

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.5530"
$ws.Range("E15").NumberFormat = "@"
$ws.Range("E15").Value = "  -0.10%  "

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "0.0₅8267"
$ws.Range("E16").NumberFormat = "@"
$ws.Range("E16").Value = "  +4.21%  "

$ws.Range("E17").NumberFormat = "@"
$ws.Range("E17").Value = "  -2.41%  "

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "26.031.30"
$ws.Range("E18").NumberFormat = "@"
$ws.Range("E18").Value = "  -0.40%  "

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "1.001"
$ws.Range("E19").NumberFormat = "@"
$ws.Range("E19").Value = "  -0.09%  "

$ws.Range("E20").NumberFormat = "@"
$ws.Range("E20").Value = "  -0.50%  "

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "187.97"
$ws.Range("E21").NumberFormat = "@"
$ws.Range("E21").Value = "  +0.93%  "

$ws.Range("E22").NumberFormat = "@"
$ws.Range("E22").Value = "  -1.30%  "

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "6.157"
$ws.Range("E23").NumberFormat = "@"
$ws.Range("E23").Value = "  -0.02%  "

$ws.Range("E24").NumberFormat = "@"
$ws.Range("E24").Value = "  -0.13%  "

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "145.28"
$ws.Range("E25").NumberFormat = "@"
$ws.Range("E25").Value = "  -2.99%  "

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "0.1217"
$ws.Range("E26").NumberFormat = "@"
$ws.Range("E26").Value = "  -2.43%  "

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "7.419"
$ws.Range("E27").NumberFormat = "@"
$ws.Range("E27").Value = "  -0.80%  "

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "15.79"
$ws.Range("E28").NumberFormat = "@"
$ws.Range("E28").Value = "  -0.61%  "

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "1.396"
$ws.Range("E29").NumberFormat = "@"
$ws.Range("E29").Value = "  +2.99%  "

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "0.05975"
$ws.Range("E30").NumberFormat = "@"
$ws.Range("E30").Value = "  -5.27%  "

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "1.254"
$ws.Range("E31").NumberFormat = "@"
$ws.Range("E31").Value = "  -1.50%  "

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "3.437"
$ws.Range("E32").NumberFormat = "@"
$ws.Range("E32").Value = "  -1.62%  "

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "3.409"
$ws.Range("E33").NumberFormat = "@"
$ws.Range("E33").Value = "  +0.01%  "

$ws.Range("E34").NumberFormat = "@"
$ws.Range("E34").Value = "  +0.45%  "

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "0.9850"
$ws.Range("E35").NumberFormat = "@"
$ws.Range("E35").Value = "  -1.09%  "

$ws.Range("E36").NumberFormat = "@"
$ws.Range("E36").Value = "  -0.54%  "

$ws.Range("E37").NumberFormat = "@"
$ws.Range("E37").Value = "  +1.20%  "

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.5682"
$ws.Range("E38").NumberFormat = "@"
$ws.Range("E38").Value = "  -5.79%  "

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.01615"
$ws.Range("E39").NumberFormat = "@"
$ws.Range("E39").Value = "  +0.13%  "

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.8508"
$ws.Range("E40").NumberFormat = "@"
$ws.Range("E40").Value = "  -2.51%  "

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "5.760"
$ws.Range("E41").NumberFormat = "@"
$ws.Range("E41").Value = "  -5.65%  "

$ws.Range("E42").NumberFormat = "@"
$ws.Range("E42").Value = "  -0.22%  "

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "1.034.72"
$ws.Range("E43").NumberFormat = "@"
$ws.Range("E43").Value = "  -6.71%  "

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "100.20"
$ws.Range("E44").NumberFormat = "@"
$ws.Range("E44").Value = "  +0.44%  "

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "1.786.12"
$ws.Range("E45").NumberFormat = "@"
$ws.Range("E45").Value = "  -1.77%  "

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.0₈107"
$ws.Range("E46").NumberFormat = "@"
$ws.Range("E46").Value = "  -1.13%  "

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "55.68"
$ws.Range("E47").NumberFormat = "@"
$ws.Range("E47").Value = "  +0.54%  "

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "1.003"
$ws.Range("E48").NumberFormat = "@"
$ws.Range("E48").Value = "  +0.32%  "

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "8.060"
$ws.Range("E49").NumberFormat = "@"
$ws.Range("E49").Value = "  +0.60%  "

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.05159"
$ws.Range("E50").NumberFormat = "@"
$ws.Range("E50").Value = "  -1.48%  "

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.4213"
$ws.Range("E51").NumberFormat = "@"
$ws.Range("E51").Value = "  -0.69%  "

